# UV-K5 Menu_RU.xlsx — "Add files via upload" commit
#
# The shared-string table had five Russian/English labels removed and
# replaced with corrected/renamed versions (placed at the end of the
# shared-string table by Excel when it rewrites the file). Each of the
# five old strings was referenced by exactly one cell in the sheet, so
# the edit is a straightforward text replacement on those five cells:
#
#   E34 : "Метод сканирования"
#         -> "Метод продолжения сканирования: TO(time out), CO(carrier out), SE(search end) "
#   E37 : "Вкл сканера листа 1"   -> "Вкл канала в лист сканирования 1"
#   E38 : "Вкл сканера листа 2"   -> "Вкл канала в лист сканирования 2"
#   E40 : "Tail tone elimination in communication through<LINE SEP>repeater"
#         -> "Tail tone elimination in communication through repeater"
#   E41 : "Чувствит микрофона"    -> "Чувствительность микрофона"
#
# The author's workbook was also left scrolled/selected at E34 (row 13
# at the top of the viewport) when it was saved; we reproduce the
# selection change that is representable through the object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E34").Value = "Метод продолжения сканирования: TO(time out), CO(carrier out), SE(search end) "
$ws.Range("E37").Value = "Вкл канала в лист сканирования 1"
$ws.Range("E38").Value = "Вкл канала в лист сканирования 2"
$ws.Range("E40").Value = "Tail tone elimination in communication through repeater"
$ws.Range("E41").Value = "Чувствительность микрофона"

# Reflect the saved cursor/selection position (cell E34) from the diff's
# updated <sheetView> (the workbook was scrolled down with E34 selected).
$ws.Activate()
$ws.Range("E34").Select()
